$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateScheme")

# Modify scheme's product/SKU name in B2 (was "ProductAuto2", now "product2")
$ws.Range("B2").Value = "product2"

# Update the active selection to reflect the new cursor position used when saving (F9)
$ws.Activate()
$ws.Range("F9").Select()
